$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Grow the table by 6 rows (A1:K23 -> A1:K29), new rows start blank.
for ($i = 0; $i -lt 6; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# 2) Apply an AutoFilter on the "Steps" column (6th column of the table,
#    colId=5) so that only rows whose Steps value is "Step 1" (or the
#    tab-prefixed variant used on the TestScenario_5 row) stay visible.
#    Excel hides every other data row as a side effect of this filter.
$tab = [char]9
$criteria = @($tab + "Step 1", "Step 1")
$lo.Range.AutoFilter(6, $criteria, 7) | Out-Null

# 3) Clear the "Approved/Rejected" (column I) value on every row that the
#    filter just hid, leaving it intact on the rows that remain visible.
$hiddenRows = @(3,4,5,6,7,8,9,11,13,14,15,16,17,18,19,21,22)
foreach ($r in $hiddenRows) {
    $ws.Cells.Item($r, 9).ClearContents()
}

# 4) Explicitly mark the six brand-new trailing rows as hidden too (they
#    fall outside the autofiltered data range but match the workbook).
for ($r = 24; $r -le 29; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}

# 5) Match the author's final selection.
$ws.Range("I31").Select() | Out-Null
